$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, derived from the authoritative diff of the
# workbook's OOXML. All target cells hold plain text (inline strings)
# in the original file, so every write below must land back in the
# worksheet as text -- including values that look numeric (e.g. "606.42",
# "1.00", "0.0590") where the exact textual formatting (trailing zeros,
# thousands dots, etc.) must be preserved exactly as authored.
$updates = @(
    @{ Cell = 'D2'; Value = '64.416.96' },
    @{ Cell = 'E2'; Value = '  +2.14%  ' },
    @{ Cell = 'D3'; Value = '2.642.12' },
    @{ Cell = 'E3'; Value = '  +0.70%  ' },
    @{ Cell = 'E4'; Value = '  +0.01%  ' },
    @{ Cell = 'D5'; Value = '606.42' },
    @{ Cell = 'E5'; Value = '  +0.20%  ' },
    @{ Cell = 'D6'; Value = '151.98' },
    @{ Cell = 'E6'; Value = '  +3.70%  ' },
    @{ Cell = 'E7'; Value = '  +0.05%  ' },
    @{ Cell = 'E8'; Value = '  +1.39%  ' },
    @{ Cell = 'D9'; Value = '0.111' },
    @{ Cell = 'E9'; Value = '  +2.48%  ' },
    @{ Cell = 'B10'; Value = 'Toncoin' },
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Cell = 'D10'; Value = '5.77' },
    @{ Cell = 'E10'; Value = '  +2.95%  ' },
    @{ Cell = 'B11'; Value = 'Cardano' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada' },
    @{ Cell = 'D11'; Value = '0.389' },
    @{ Cell = 'E11'; Value = '  +7.34%  ' },
    @{ Cell = 'E12'; Value = '  -0.51%  ' },
    @{ Cell = 'D13'; Value = '27.86' },
    @{ Cell = 'E13'; Value = '  +2.35%  ' },
    @{ Cell = 'D14'; Value = '3.118.89' },
    @{ Cell = 'E14'; Value = '  +0.75%  ' },
    @{ Cell = 'D15'; Value = '64.274.82' },
    @{ Cell = 'E15'; Value = '  +2.13%  ' },
    @{ Cell = 'E16'; Value = '  +4.03%  ' },
    @{ Cell = 'D17'; Value = '2.735.53' },
    @{ Cell = 'E17'; Value = '  +4.39%  ' },
    @{ Cell = 'D18'; Value = '12.25' },
    @{ Cell = 'E18'; Value = '  +8.58%  ' },
    @{ Cell = 'D19'; Value = '4.68' },
    @{ Cell = 'E19'; Value = '  +4.21%  ' },
    @{ Cell = 'D20'; Value = '353.29' },
    @{ Cell = 'E20'; Value = '  +4.23%  ' },
    @{ Cell = 'D21'; Value = '6.99' },
    @{ Cell = 'E21'; Value = '  +1.77%  ' },
    @{ Cell = 'E22'; Value = '  +0.43%  ' },
    @{ Cell = 'D23'; Value = '5.74' },
    @{ Cell = 'E23'; Value = '  +3.21%  ' },
    @{ Cell = 'D24'; Value = '66.89' },
    @{ Cell = 'E24'; Value = '  +0.52%  ' },
    @{ Cell = 'E25'; Value = '  +13.55%  ' },
    @{ Cell = 'D26'; Value = '1.72' },
    @{ Cell = 'E26'; Value = '  +5.93%  ' },
    @{ Cell = 'D27'; Value = '9.35' },
    @{ Cell = 'E27'; Value = '  +7.56%  ' },
    @{ Cell = 'D28'; Value = '8.26' },
    @{ Cell = 'E28'; Value = '  +4.20%  ' },
    @{ Cell = 'E29'; Value = '  +2.26%  ' },
    @{ Cell = 'D30'; Value = '550.65' },
    @{ Cell = 'E30'; Value = '  +2.52%  ' },
    @{ Cell = 'E31'; Value = '  +0.05%  ' },
    @{ Cell = 'E32'; Value = '  +1.96%  ' },
    @{ Cell = 'E33'; Value = '  +8.89%  ' },
    @{ Cell = 'D34'; Value = '1.80' },
    @{ Cell = 'E34'; Value = '  +2.62%  ' },
    @{ Cell = 'D35'; Value = '5.35' },
    @{ Cell = 'E35'; Value = '  +1.62%  ' },
    @{ Cell = 'D36'; Value = '167.69' },
    @{ Cell = 'E36'; Value = '  -0.75%  ' },
    @{ Cell = 'D37'; Value = '2.03' },
    @{ Cell = 'E37'; Value = '  +8.80%  ' },
    @{ Cell = 'D38'; Value = '0.413' },
    @{ Cell = 'E38'; Value = '  +2.41%  ' },
    @{ Cell = 'D39'; Value = '1.00' },
    @{ Cell = 'E39'; Value = '  -0.03%  ' },
    @{ Cell = 'E40'; Value = '  +3.41%  ' },
    @{ Cell = 'E41'; Value = '  +0.05%  ' },
    @{ Cell = 'D42'; Value = '169.47' },
    @{ Cell = 'E42'; Value = '  +0.62%  ' },
    @{ Cell = 'D43'; Value = '40.26' },
    @{ Cell = 'E43'; Value = '  +1.36%  ' },
    @{ Cell = 'D44'; Value = '3.96' },
    @{ Cell = 'E44'; Value = '  +5.50%  ' },
    @{ Cell = 'D45'; Value = '0.0590' },
    @{ Cell = 'E45'; Value = '  +4.03%  ' },
    @{ Cell = 'D46'; Value = '21.85' },
    @{ Cell = 'E46'; Value = '  -1.68%  ' },
    @{ Cell = 'D47'; Value = '0.633' },
    @{ Cell = 'E47'; Value = '  +1.61%  ' },
    @{ Cell = 'D48'; Value = '2.04' },
    @{ Cell = 'E48'; Value = '  +15.52%  ' },
    @{ Cell = 'E49'; Value = '  +2.84%  ' },
    @{ Cell = 'D50'; Value = '0.0969' },
    @{ Cell = 'E50'; Value = '  +1.06%  ' },
    @{ Cell = 'D51'; Value = '19.50' },
    @{ Cell = 'E51'; Value = '  +5.53%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)

    # Column D sometimes contains values that Excel would otherwise
    # auto-convert to a Number (e.g. "606.42" -> 606.42, dropping the
    # "1.00"/"0.0590"-style trailing zeros). Force the cell to Text so
    # the literal string is preserved, then restore the original
    # (default/"Normal") cell style so no visible formatting changes
    # leak into the sheet.
    if ($u.Cell.StartsWith("D")) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
